# Fix: combustion/biomass filters were applied on the wrong (unfiltered)
# dataframes; correcting downstream aggregated numbers per plant sheet.

$wb = $excel.ActiveWorkbook

# --- Strausberg ---
$ws = $wb.Worksheets.Item("Strausberg")
$ws.Cells.Item(5, 3).Value = 10299.7
$ws.Cells.Item(5, 4).Value = 12044.2
$ws.Cells.Item(5, 5).Value = 9

# --- Ruedersdorf bei Berlin ---
$ws = $wb.Worksheets.Item("Rüdersdorf bei Berlin")
$ws.Cells.Item(4, 3).Value = 1
$ws.Cells.Item(4, 4).Value = 0
$ws.Cells.Item(4, 5).Value = 0

# --- Gruenheide (Mark) ---
$ws = $wb.Worksheets.Item("Grünheide (Mark)")
$ws.Cells.Item(2, 3).Value = 194.7
$ws.Cells.Item(2, 4).Value = 328.6
$ws.Cells.Item(2, 5).Value = 4

# --- Bocholt ---
$ws = $wb.Worksheets.Item("Bocholt")
$ws.Cells.Item(2, 3).Value = 8.75
$ws.Cells.Item(2, 4).Value = 7.9
$ws.Cells.Item(2, 5).Value = 8

$ws.Cells.Item(7, 3).Value = 1
$ws.Cells.Item(7, 4).Value = 5.82
$ws.Cells.Item(7, 5).Value = 1

$ws.Cells.Item(10, 3).Value = 2925.8
$ws.Cells.Item(10, 4).Value = 4657.2
$ws.Cells.Item(10, 5).Value = 51

$ws.Cells.Item(12, 3).Value = 1039.6
$ws.Cells.Item(12, 4).Value = 79
$ws.Cells.Item(12, 5).Value = 4

# --- Zwickau ---
$ws = $wb.Worksheets.Item("Zwickau")
$ws.Cells.Item(8, 3).Value = 20314.1
$ws.Cells.Item(8, 4).Value = 22719.3
$ws.Cells.Item(8, 5).Value = 35

# --- Ingolstadt ---
$ws = $wb.Worksheets.Item("Ingolstadt")
$ws.Cells.Item(7, 3).Value = 2.1
$ws.Cells.Item(7, 4).Value = 29
$ws.Cells.Item(7, 5).Value = 2

$ws.Cells.Item(9, 3).Value = 2348.4
$ws.Cells.Item(9, 4).Value = 3541.82
$ws.Cells.Item(9, 5).Value = 51

$ws.Cells.Item(10, 3).Value = 230.3
$ws.Cells.Item(10, 5).Value = 1

# --- Kassel ---
$ws = $wb.Worksheets.Item("Kassel")
$ws.Cells.Item(4, 3).Value = 42.5
$ws.Cells.Item(4, 4).Value = 83
$ws.Cells.Item(4, 5).Value = 4

# Row 9 ("Kondensationsmaschine mit Entnahme" / "Rohbraunkohlen") was an
# erroneous leftover row from the unfiltered dataframe; remove it so every
# row below shifts up by one.
$ws.Rows.Item(9).Delete()

# The row that is now row 13 (previously row 14) also needs corrected
# aggregate values.
$ws.Cells.Item(13, 3).Value = 7307.77
$ws.Cells.Item(13, 4).Value = 11801.3
$ws.Cells.Item(13, 5).Value = 77

# --- Kiel ---
$ws = $wb.Worksheets.Item("Kiel")
$ws.Cells.Item(7, 3).Value = 13
$ws.Cells.Item(7, 4).Value = 31.9
$ws.Cells.Item(7, 5).Value = 2

$ws.Cells.Item(10, 3).Value = 194611.9
$ws.Cells.Item(10, 4).Value = 4038804.2
$ws.Cells.Item(10, 5).Value = 83

$ws.Cells.Item(12, 3).Value = 1565
$ws.Cells.Item(12, 4).Value = 0
$ws.Cells.Item(12, 5).Value = 0
